$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.553.24"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "2.761.19"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.58"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").Value = "3.194.36"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").Value = "2.767.26"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.892"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "51.591.77"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.49"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.63"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.02"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0823"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.98"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("E40").Value = "  +9.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.33%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.114"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("E45").Value = "  +14.05%  "
$ws.Range("D46").Value = "2.090.73"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("E49").Value = "  +5.65%  "
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.24%  "
